$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "gender"

# Select column F (as the last user interaction before save)
$ws.Columns("F").Select() | Out-Null
